$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "SC 92" row (originally row 28) first so the "RM 232" row
# index (26) stays valid while we still need it.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()
